$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.386.36'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '1.871.80'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.52'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4709'
$ws.Range("E7").Value = '  -0.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2867'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06492'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.83'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '100.58'
$ws.Range("E11").Value = '  +3.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07799'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '1.872.60'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7280'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.167'
$ws.Range("E15").Value = '  -1.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.00'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '30.375.13'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007487'
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").Value = '2.115.79'
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.325'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.336'
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.042'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.94'
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09678'
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.322'
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.491'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.222'
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04804'
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6899'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.712'
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01896'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.841'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.20'
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.302'
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.951'
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4209'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8242'
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.79'
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.712'
$ws.Range("E47").Value = '  +2.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.005'
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05757'
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '883.52'
$ws.Range("E51").Value = '  -3.69%  '
